$d = $word.ActiveDocument

# 1. Delete the whole paragraph "- opcje i timer na dwóch wjeżdżających stronach"
#    (paragraph 7), including its paragraph mark.
[void]$d.Paragraphs.Item(7).Range.Delete()

# 2. Delete the whole paragraphs "- zamiana sekund na minuty" and
#    "- wyświetlanie liczby pozostałych rund" (now paragraphs 8 and 9,
#    after the previous deletion). Delete from the end first so indices
#    of not-yet-deleted paragraphs stay stable.
[void]$d.Paragraphs.Item(9).Range.Delete()
[void]$d.Paragraphs.Item(8).Range.Delete()

# 3. Relocate the "_GoBack" bookmark from the end of "Przycisk do
#    potwierdzenia opcji" (paragraph 6) to the end of "- losowanie czasu
#    wydawania dźwięków" (now paragraph 7), right before its paragraph
#    mark, matching the target markup exactly:
#      ...<w:t>losowanie czasu wydawania dźwięków</w:t></w:r>
#      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
#      <w:bookmarkEnd w:id="0"/>
#
#    A collapsed Range sitting exactly at "paragraph end minus one"
#    (i.e. immediately before the pilcrow) cannot be passed straight to
#    Bookmarks.Add in this host, so a trailing placeholder character is
#    inserted first to move that boundary out of the way, the bookmark
#    is added at the (now safe) position just before the placeholder,
#    and the placeholder is removed again. Bookmarks track by position,
#    so the bookmark stays put once the placeholder goes away.
$target = $d.Paragraphs.Item(7)
$insertPos = $target.Range.End - 1
$placeholderRange = $d.Range($insertPos, $insertPos)
$placeholderRange.InsertAfter("X")

$target = $d.Paragraphs.Item(7)
$bmPos = $target.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Bookmarks.Item("_GoBack").Range.End
$placeholderRange = $d.Range($placeholder, $placeholder + 1)
$placeholderRange.Delete()
